# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and push
#    the existing 2022-Q2 / 2022-Q1 rows down by one.
# 2) Insert a brand-new "2022-Q3" worksheet (with the per-fund breakdown)
#    right before the existing "2022-Q2" worksheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Update the summary sheet -------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push row 2 (2022-Q2) down to row 3, and row 3 (2022-Q1) down to row 4
# first, so nothing gets clobbered while we still need to read it.
# A4 is a brand new cell, so pick up the index-column style (s="2") that
# A2/A3 already carry before writing its value.
$summary.Range("A2").Copy($summary.Range("A4"))
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 3.61

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 10
$summary.Range("D3").Value = 0.66

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.42

# --- 2. Create the new "2022-Q3" worksheet ----------------------------------
# Worksheets.Add() on its own produces a worksheet whose cells silently
# refuse copy/format propagation in this host, so instead duplicate the
# existing "2022-Q2" sheet (which already carries the right header / index
# column styles) and overwrite its contents in place.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The source sheet has 10 data rows (2..11); the new sheet only needs 5
# (2..6), so drop the extra rows entirely.
$q3.Range("A7:H11").Clear()

# Header row (unchanged text, but re-assert to be safe).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'159883"
$q3.Range("C2").Value = "永赢中证全指医疗器械ETF"
$q3.Range("D2").Value = "'10.67"
$q3.Range("E2").Value = "'98.26"
$q3.Range("F2").Value = "'2.65"
$q3.Range("G2").Value = "'0.2828"
$q3.Range("H2").Value = 9

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'515860"
$q3.Range("C3").Value = "嘉实新兴科技100ETF"
$q3.Range("D3").Value = "'1.92"
$q3.Range("E3").Value = "'98.86"
$q3.Range("F3").Value = "'4.13"
$q3.Range("G3").Value = "'0.0793"
$q3.Range("H3").Value = 5

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'159898"
$q3.Range("C4").Value = "招商中证全指医疗器械ETF"
$q3.Range("D4").Value = "'1.35"
$q3.Range("E4").Value = "'99.00"
$q3.Range("F4").Value = "'2.56"
$q3.Range("G4").Value = "'0.0346"
$q3.Range("H4").Value = 10

# Row 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'159797"
$q3.Range("C5").Value = "汇添富中证全指医疗器械ETF"
$q3.Range("D5").Value = "'0.54"
$q3.Range("E5").Value = "'94.87"
$q3.Range("F5").Value = "'2.63"
$q3.Range("G5").Value = "'0.0142"
$q3.Range("H5").Value = 10

# Row 6
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'515870"
$q3.Range("C6").Value = "嘉实中证先进制造100策略ETF"
$q3.Range("D6").Value = "'0.36"
$q3.Range("E6").Value = "'98.05"
$q3.Range("F6").Value = "'2.58"
$q3.Range("G6").Value = "'0.0093"
$q3.Range("H6").Value = 8

# The leading apostrophes above force text storage (matching the other
# quarters' sheets, where these numeric-looking columns are strings); strip
# the resulting quote-prefix styling so the cells fall back to the default
# (unstyled) format like their siblings.
$q3.Range("B2:B6").Style = "Normal"
$q3.Range("D2:G6").Style = "Normal"

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("2022-Q1").Activate()
